$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Mã Barcode" column (C) and its row-2 value (5555551412) are being removed
# from the "quy cach" (spec) check export - delete the entire column, which
# shifts Tên sản phẩm / Quy cách / Trạng thái left by one.
$ws.Columns(3).EntireColumn.Delete()

# Reflect the operator's follow-up click into the (now empty) F column.
[void]$ws.Range("F9").Select()
